$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for the name / phone number columns (rows 2-10)
$names = @("Ana Carolina", "Anna Lena", "Jhonatan", "Julia", "Brenda", "Caeh", "Laura", "Juliana", "Larissa")
$phones = @("5532984884745", "5532988191429", "5532988923958", "5532984427134", "5532988213508", "5532991418096", "5532984700198", "5532988858637", "5532999798051")

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).NumberFormat = "General"
    $ws.Cells.Item($row, 2).Value = $phones[$i]
}

# Adjust column widths to match the new content
$ws.Columns.Item(1).ColumnWidth = 16.81640625
$ws.Columns.Item(3).ColumnWidth = 34.26953125

# Update the active selection as seen in the final workbook
$ws.Range("C2:C10").Select()
